$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) { return $candidate }
    }
    return $null
}

# --- Shape "TextBox 47" (id=48): "[command commits address book]"
#     -> "[command commits wish book]", split into
#     "command " / "commits wish " / "book]" runs.
$sh48 = Get-ShapeById $s 48
$tr48 = $sh48.TextFrame.TextRange
$mid48 = $tr48.Characters(10, 16)
$mid48.Text = "commits wish "
# Restore the autofit-computed box height back to its original value
# (text editing nudges it by a handful of EMUs through layout rounding).
$sh48.Height = 50.91237

# --- Shape "Rectangle: Rounded Corners 50" (id=51):
#     "Purge redundant states and then save address book to addressBookStateList "
#     -> "Purge redundant states and then save wish book to wishBookStateList "
$sh51 = Get-ShapeById $s 51
$tr51 = $sh51.TextFrame.TextRange
$run1_51 = $tr51.Characters(1, 53)
$run1_51.Text = "Purge redundant states and then save wish book to "
$run2_51 = $tr51.Characters(51, 20)
$run2_51.Text = "wishBookStateList"
